# Applies the "Text Functions" Data-sheet enhancements:
#   - J: Full Address = TEXTJOIN(" ", TRUE, Street, CityStateZip)
#   - L: First Name   = LEFT(FullName, FIND(" ", FullName) - 1)
#   - M: Last Name    = RIGHT(FullName, LEN(FullName) - FIND(" ", FullName))
#   - O: Applicant ID (Short) = RIGHT(ApplicantID, 7)   (filled one row further -> row 22)
#   - T/U/V: second "Find and MID (State)" helper block (same idea as Q/R/S,
#     just without the extra +LEN(", ") offset, so it captures the leading space)
#
# Matches the pattern already used for columns Q/R/S in this sheet: the first
# data row (row 2) gets its own standalone formula, and rows 3-21 are filled
# as one shared-formula block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- J: Full Address ---------------------------------------------------
$ws.Range("J2").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,F2,G2)'
$ws.Range("J3:J21").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,F3,G3)'

# --- L: First Name -------------------------------------------------------
$ws.Range("L2").Formula = '=LEFT(B2,FIND(" ",B2)-1)'
$ws.Range("L3:L21").Formula = '=LEFT(B3,FIND(" ",B3)-1)'

# --- M: Last Name ----------------------------------------------------------
$ws.Range("M2").Formula = '=RIGHT(B2,LEN(B2)-FIND(" ",B2))'
$ws.Range("M3:M21").Formula = '=RIGHT(B3,LEN(B3)-FIND(" ",B3))'

# --- O: Applicant ID (Short) -- filled one extra row (through row 22) -----
$ws.Range("O2").Formula = '=RIGHT(A2,7)'
$ws.Range("O3:O22").Formula = '=RIGHT(A3,7)'

# --- T/U/V: second Find-and-MID (State) helper block ------------------
$ws.Range("T2").Formula = '=FIND(",",G2)+1'
$ws.Range("T3:T21").Formula = '=FIND(",",G3)+1'

$ws.Range("U2").Formula = '=FIND(",",G2,T2)'
$ws.Range("U3:U21").Formula = '=FIND(",",G3,T3)'

$ws.Range("V2").Formula = '=MID(G2,T2,U2-T2)'
$ws.Range("V3:V21").Formula = '=MID(G3,T3,U3-T3)'

Write-Output "done"
